$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update B2: CasesTab query (Age no longer wrapped in COALESCE/CAST; Cohort now joins df_cohort) ---
$b2Query = @"
SELECT
    DISTINCT c.case_id AS "Case ID",
    st.clinical_study_designation AS "Study Code",
    st.clinical_study_type AS "Study Type",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis",
    diag.stage_of_disease AS "Stage Of Disease",
    dmg.patient_age_at_enrollment AS "Age", 
    COALESCE(TRIM(dmg.sex), '') AS "Sex",   
    COALESCE(dmg.neutered_indicator, '') AS "Neutered Status",  
    COALESCE(
        CASE 
            WHEN dmg.weight = CAST(dmg.weight AS INT) THEN CAST(dmg.weight AS INT)
            ELSE dmg.weight 
        END, 
    '') AS "Weight (kg)",        
    COALESCE(diag.best_response, '') AS "Response to Treatment",  
    COALESCE(coh.cohort_description, '') AS "Cohort"  
FROM 
    df_case c
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_demographic dmg ON dmg."case.case_id" = c.case_id
JOIN 
    df_diagnosis diag ON diag."case.case_id" = c.case_id
JOIN 
    df_enrollment enr ON enr."case.case_id" = c.case_id
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_sample smp ON smp."case.case_id" = c.case_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_case_file cf ON cf."sample.sample_id" = smp.sample_id
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN
    df_cohort coh ON coh."study.clinical_study_designation" = st.clinical_study_designation
WHERE
     st.clinical_study_designation = 'MGT01' AND dmg.breed = 'Labrador Retriever'
ORDER BY 
    c.case_id ASC
LIMIT 100;
"@
$ws.Range("B2").Value = $b2Query

# --- Update B4: CaseFilesTab query (original_file_size -> file_size) ---
$b4Query = @"
SELECT DISTINCT
    cf.file_name AS "File Name",
    CASE
        WHEN cf.file_name LIKE '%.bai' THEN 'bai'
        WHEN cf.file_name LIKE '%.bam' THEN 'bam'
        WHEN cf.file_name LIKE '%.csv' THEN 'csv'
        WHEN cf.file_name LIKE '%.doc' THEN 'doc'
        WHEN cf.file_name LIKE '%.docx' THEN 'docx'
        WHEN cf.file_name LIKE '%.gz' THEN 'gz'
        WHEN cf.file_name LIKE '%.pdf' THEN 'pdf'
        WHEN cf.file_name LIKE '%.rtf' THEN 'rtf'
        WHEN cf.file_name LIKE '%.tbi' THEN 'tbi'
        WHEN cf.file_name LIKE '%.tif' THEN 'tif'
        WHEN cf.file_name LIKE '%.xls' THEN 'xls'
        WHEN cf.file_name LIKE '%.xlsx' THEN 'xlsx'
        ELSE 'Unknown'
    END AS "Format",
    cf.file_type AS "File Type",
    CASE     
    WHEN cf.file_size >= 1024 * 1024 * 1024 THEN 
        ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB' 
    WHEN cf.file_size >= 1024 * 1024 THEN 
        ROUND(cf.file_size / (1024.0 * 1024.0), 2) || ' MB' 
    WHEN cf.file_size >= 1024 THEN 
        ROUND(cf.file_size / 1024.0, 2) || ' KB' 
    ELSE 
        ROUND(cf.file_size, 2) || ' Bytes' 
END AS "Size",
    'sample' AS "Association",
    cf.file_description AS "Description",
    smp.sample_id AS "Sample ID",
    c.case_id AS "Case ID",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis"
FROM 
    df_case_file cf
JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
JOIN 
    df_case c ON smp."case.case_id" = c.case_id
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_demographic dmg ON dmg."case.case_id" = c.case_id
JOIN 
    df_diagnosis diag ON diag."case.case_id" = c.case_id
JOIN 
    df_enrollment enr ON enr."case.case_id" = c.case_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = 'MGT01' AND dmg.breed = 'Labrador Retriever'
ORDER BY 
    sf.file_name ASC
LIMIT 100;
"@
$ws.Range("B4").Value = $b4Query

# Preserve the (Excel max-capped) auto-fit row height of row 4 after the text change
$ws.Rows.Item(4).Height = 409.6

# --- Re-apply formatting (12pt Calibri + wrap text) on the edited cells, matching the rest of the sheet ---
foreach ($addr in @("B2", "C2", "B4")) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 12
    $rng.WrapText = $true
}

# --- Update the active view: scrolled to show row 3 at top, with B4 selected ---
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1

$wb.Save()
